$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 9945
$ws.Range("C3").Value = 110.77
$ws.Range("D3").Value = 110.16
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = -0.55
$ws.Range("G3").Value = 42608.639097222222
